$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated coin data (price & 1h volume%), and for row 51 a coin swap
# (Decentraland -> NEARProtocol) as captured by the scheduled GitHub Actions scrape.

$ws.Range("D2").Value = "'30.654.11"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").Value = "'1.892.55"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.52%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'240.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.16%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.18%  "
$ws.Range("D7").Value = "'0.4905"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.00%  "
$ws.Range("D8").Value = "'0.2945"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.02%  "
$ws.Range("D9").Value = "'0.06708"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.54%  "
$ws.Range("D10").Value = "'1.923.13"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.14%  "
$ws.Range("D11").Value = "'17.20"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.94%  "
$ws.Range("D12").Value = "'0.07351"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'5.148"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.11%  "
$ws.Range("D14").Value = "'88.31"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.01%  "
$ws.Range("D15").Value = "'0.6684"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.71%  "
$ws.Range("D16").Value = "'30.591.99"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.15%  "
$ws.Range("D17").Value = "'0.000007875"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.66%  "
$ws.Range("D18").Value = "'13.42"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.53%  "
$ws.Range("D19").Value = "'1.002"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.24%  "
$ws.Range("D20").Value = "'2.144.90"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.49%  "
$ws.Range("D21").Value = "'5.330"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +12.29%  "
$ws.Range("E22").Value = "  -0.21%  "
$ws.Range("D23").Value = "'190.05"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.99%  "
$ws.Range("D24").Value = "'6.201"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.94%  "
$ws.Range("D25").Value = "'9.525"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.19%  "
$ws.Range("D26").Value = "'161.63"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.17%  "
$ws.Range("D27").Value = "'18.45"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.27%  "
$ws.Range("D28").Value = "'1.934"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.86%  "
$ws.Range("D29").Value = "'1.464"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.45%  "
$ws.Range("D30").Value = "'4.402"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.81%  "
$ws.Range("D31").Value = "'0.09156"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.90%  "
$ws.Range("D32").Value = "'4.056"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.78%  "
$ws.Range("D33").Value = "'0.05247"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.17%  "
$ws.Range("D34").Value = "'0.7434"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.85%  "
$ws.Range("D35").Value = "'1.101"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.39%  "
$ws.Range("D36").Value = "'2.722"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.13%  "
$ws.Range("D37").Value = "'0.01826"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.33%  "
$ws.Range("D38").Value = "'2.688"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.82%  "
$ws.Range("D39").Value = "'0.9153"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.05%  "
$ws.Range("D40").Value = "'2.069"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.65%  "
$ws.Range("D41").Value = "'74.92"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +31.61%  "
$ws.Range("D42").Value = "'0.4429"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.23%  "
$ws.Range("D43").Value = "'5.923"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.96%  "
$ws.Range("D44").Value = "'106.20"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.19%  "
$ws.Range("D45").Value = "'0.9928"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.60%  "
$ws.Range("D46").Value = "'0.1379"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.84%  "
$ws.Range("D47").Value = "'7.563"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.56%  "
$ws.Range("D48").Value = "'35.44"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.38%  "
$ws.Range("D49").Value = "'9.037"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.42%  "
$ws.Range("D50").Value = "'0.05837"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.19%  "
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "'1.426"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.71%  "
